# "complete block 3 notebook"
# Slide 8 ("Coding") content placeholder loses its first bullet
# ("Load dataset") and the final bullet ("Try different input features")
# gets re-touched (picks up dirty="0") as part of the same edit.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(8)
$sh = $s.Shapes.Item(2)

# Rewrite the body text without the "Load dataset" bullet; PowerPoint
# regenerates the runs for every remaining paragraph, so the trailing
# bullet also ends up freshly marked dirty="0".
$sh.TextFrame.TextRange.Text = "Z-score inputs`rOptimize model with gradient descent`rTry different input features"
